$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mandarin")
$ws.Activate()

# Replace header row with trial-structure column names (holdover testing
# data from October is being cleared out, headers reflect the new layout)
$ws.Range("B1").Value = "trial1"
$ws.Range("C1").Value = "trial2"
$ws.Range("D1").Value = "question"
$ws.Range("A1").Value = "trialNum"

# Remove the now-unused overflow rows (8-12) left over from the October
# testing data
$ws.Range("A8:D12").EntireRow.Delete()
$ws.Range("A4:A7").Formula = "=A3+1"

$ws.Range("D17").Select()
